# Updated symbol list on Tue Jan 10 07:21:43 UTC 2023 with GitHub Actions
# This script updates Price (D) and Volume(1h) (E) columns for the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the target cells keep their original text/"@" format so that values such as
# "274.62" or "-1.10%" are stored as literal text, matching the source data format,
# rather than being auto-converted to numbers/percentages by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "274.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.10%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.08%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.877"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.37%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06326"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.30%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.874"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.52%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.56%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.256"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "33.37%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8687"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.35%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1538"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.08%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05002"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.99%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07479"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.68%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02943"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-6.88%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09014"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.44%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.32%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006329"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.89%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005792"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.55%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.447"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.18%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3127"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.06%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1334"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.16%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.902"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.40%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04349"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.51%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001178"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.08%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004244"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.43%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.02%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.28%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.34%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006725"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.34%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1165"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.13%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.64%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01070"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.95%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005307"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.35%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-33.05%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.490"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-37.27%"
